$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

$ws.Range("A2").Value = 2000
$ws.Range("A3").Value = 2001
$ws.Range("A4").Value = 2002
$ws.Range("A5").Value = 2100
$ws.Range("A6").Value = 2101
$ws.Range("A7").Value = 2102
$ws.Range("A8").Value = 2103
$ws.Range("A9").Value = 2104
$ws.Range("A10").Value = 2105
$ws.Range("A11").Value = 2106
$ws.Range("A12").Value = 2107
$ws.Range("A13").Value = 2108
$ws.Range("A14").Value = 2109
$ws.Range("A15").Value = 2110

$ws.Range("L27").Select()
